# Refresh the NATMI LR-pair table (Artn -> Gfra1) with newly recomputed TPM
# values. The sending/target cluster set grows from {FAPs, MuSCs} to
# {FAPs, Inflammatory-Mac, MuSCs}, so the 3x3 grid of rows (2..7) becomes a
# 3x4 grid (rows 2..13): every sending cluster now also pairs with the new
# "Resolving-Mac" target cluster, and "Inflammatory-Mac" is an entirely new
# sending cluster with its own four target rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Artn", "Gfra1", "ECs", 3, 1, 1.415594, 4.246782, 0.7538443241324221, 0.7538443241324221, 2, 0.6666666666666666, 0.5471986666666667, 1.641596, 0.02624837670119588, 0.02624837670119588, 0.7746111493413333, 6.971500344072, 0.01978718979388622, 0.01978718979388622),
    @("FAPs", "Artn", "Gfra1", "FAPs", 3, 1, 1.415594, 4.246782, 0.7538443241324221, 0.7538443241324221, 3, 1, 17.676258, 53.028774, 0.8479060840514853, 0.8479060840514852, 25.022404767252, 225.201642905268, 0.6391891888595606, 0.6391891888595604),
    @("FAPs", "Artn", "Gfra1", "MuSCs", 3, 1, 1.415594, 4.246782, 0.7538443241324221, 0.7538443241324221, 3, 1, 2.604410666666667, 7.813232000000001, 0.1249300417336775, 0.1249300417336775, 3.686788113269333, 33.181093019424, 0.09417780287455943, 0.09417780287455943),
    @("FAPs", "Artn", "Gfra1", "Resolving-Mac", 3, 1, 1.415594, 4.246782, 0.7538443241324221, 0.7538443241324221, 1, 0.3333333333333333, 0.01908533333333333, 0.057256, 0.0009154975136414022, 0.0009154975136414022, 0.02701708335466666, 0.243153750192, 0.0006901426044159157, 0.0006901426044159157),
    @("Inflammatory-Mac", "Artn", "Gfra1", "ECs", 1, 0.3333333333333333, 0.2356576666666667, 0.706973, 0.1254944528268394, 0.1254944528268394, 2, 0.6666666666666666, 0.5471986666666667, 1.641596, 0.02624837670119588, 0.02624837670119588, 0.1289515609897778, 1.160564048908, 0.003294025671709338, 0.003294025671709338),
    @("Inflammatory-Mac", "Artn", "Gfra1", "FAPs", 1, 0.3333333333333333, 0.2356576666666667, 0.706973, 0.1254944528268394, 0.1254944528268394, 3, 1, 17.676258, 53.028774, 0.8479060840514853, 0.8479060840514852, 4.165545715678, 37.489911441102, 0.1064075100665893, 0.1064075100665893),
    @("Inflammatory-Mac", "Artn", "Gfra1", "MuSCs", 1, 0.3333333333333333, 0.2356576666666667, 0.706973, 0.1254944528268394, 0.1254944528268394, 3, 1, 2.604410666666667, 7.813232000000001, 0.1249300417336775, 0.1249300417336775, 0.6137493407484444, 5.523744066736, 0.01567802722900208, 0.01567802722900208),
    @("Inflammatory-Mac", "Artn", "Gfra1", "Resolving-Mac", 1, 0.3333333333333333, 0.2356576666666667, 0.706973, 0.1254944528268394, 0.1254944528268394, 1, 0.3333333333333333, 0.01908533333333333, 0.057256, 0.0009154975136414022, 0.0009154975136414022, 0.004497605120888888, 0.040478446088, 0.0001148898595387598, 0.0001148898595387598),
    @("MuSCs", "Artn", "Gfra1", "ECs", 2, 0.6666666666666666, 0.2265816666666667, 0.679745, 0.1206612230407385, 0.1206612230407385, 2, 0.6666666666666666, 0.5471986666666667, 1.641596, 0.02624837670119588, 0.02624837670119588, 0.1239851858911111, 1.11586667302, 0.003167161235600319, 0.003167161235600319),
    @("MuSCs", "Artn", "Gfra1", "FAPs", 2, 0.6666666666666666, 0.2265816666666667, 0.679745, 0.1206612230407385, 0.1206612230407385, 3, 1, 17.676258, 53.028774, 0.8479060840514853, 0.8479060840514852, 4.005115998070001, 36.04604398263, 0.1023093851253354, 0.1023093851253354),
    @("MuSCs", "Artn", "Gfra1", "MuSCs", 2, 0.6666666666666666, 0.2265816666666667, 0.679745, 0.1206612230407385, 0.1206612230407385, 3, 1, 2.604410666666667, 7.813232000000001, 0.1249300417336775, 0.1249300417336775, 0.5901117095377779, 5.311005385840001, 0.01507421163011603, 0.01507421163011603),
    @("MuSCs", "Artn", "Gfra1", "Resolving-Mac", 2, 0.6666666666666666, 0.2265816666666667, 0.679745, 0.1206612230407385, 0.1206612230407385, 1, 0.3333333333333333, 0.01908533333333333, 0.057256, 0.0009154975136414022, 0.0009154975136414022, 0.004324386635555556, 0.03891947972, 0.0001104650496867267, 0.0001104650496867267)
)

$startRow = 2

# First pass: write column A (Sending cluster) for every row so any new
# cluster label introduced there is interned before later columns are touched.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
}

# Second pass: write the remaining columns (B..T) for every row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
